$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.922680974006653
$ws.Range("B1").Value = 2.028523921966553
$ws.Range("C1").Value = 2.071424722671509
$ws.Range("D1").Value = 2.535306453704834
$ws.Range("E1").Value = 3.659120321273804
